# Added isNotAutoSelect into lowcode.Spin
# This reorders the data rows (rows 2-24) of Sheet1 to reflect the new
# row order produced by the upstream change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(901,16,15,45,60,60),
    @(301,6,45,30,60,45),
    @(801,3,67,65,52,45),
    @(401,9,48,67,75,45),
    @(701,3,90,45,97,15),
    @(201,9,30,15,45,30),
    @(1201,2,10,10,10,10),
    @(101,9,30,15,60,15),
    @(902,1,0,0,0,0),
    @(601,9,60,67,60,42),
    @(1202,2,10,10,10,10),
    @(1203,3,15,15,15,15),
    @(1001,18,30,75,60,72),
    @(501,9,52,30,75,45),
    @(2,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(802,0,4,5,4,0),
    @(502,0,4,0,0,0),
    @(1,0,2,2,2,2),
    @(1101,0,15,30,30,0),
    @(602,0,0,4,0,9),
    @(402,0,0,4,0,0),
    @(702,0,0,0,4,0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$col - 1]
    }
}
